$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the text-format number-format on D-column cells whose new value would
# otherwise be auto-parsed by Excel as a numeric literal (losing the exact
# textual representation, e.g. trailing zeros like "1.00" or "166.05").
$textCells = @("D5", "D6", "D11", "D12", "D14", "D19", "D21", "D25", "D28", "D32", "D34", "D35", "D36", "D38", "D40", "D43", "D48", "D49")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Apply the updated price (column D) / volume-change (column E) values.
$ws.Range("D2").Value = "69.333.04"
$ws.Range("E2").Value = "  +2.35%  "
$ws.Range("D3").Value = "3.781.87"
$ws.Range("E3").Value = "  +0.28%  "
$ws.Range("E4").Value = "  -0.48%  "
$ws.Range("D5").Value = "624.75"
$ws.Range("D6").Value = "166.05"
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("D7").Value = "3.781.41"
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("E9").Value = "  +1.30%  "
$ws.Range("E10").Value = "  +2.84%  "
$ws.Range("D11").Value = "0.454"
$ws.Range("E11").Value = "  +2.12%  "
$ws.Range("D12").Value = "6.68"
$ws.Range("E12").Value = "  +1.03%  "
$ws.Range("D14").Value = "35.66"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("D15").Value = "4.420.51"
$ws.Range("E15").Value = "  +0.39%  "
$ws.Range("D16").Value = "3.792.52"
$ws.Range("E16").Value = "  +0.43%  "
$ws.Range("D17").Value = "69.386.27"
$ws.Range("E17").Value = "  +2.41%  "
$ws.Range("E18").Value = "  -2.74%  "
$ws.Range("D19").Value = "7.12"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "468.20"
$ws.Range("E21").Value = "  +2.42%  "
$ws.Range("E22").Value = "  +1.63%  "
$ws.Range("E23").Value = "  +2.04%  "
$ws.Range("E24").Value = "  +4.61%  "
$ws.Range("D25").Value = "83.25"
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  +1.74%  "
$ws.Range("E27").Value = "  +3.78%  "
$ws.Range("D28").Value = "10.02"
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "3.934.81"
$ws.Range("E30").Value = "  +0.42%  "
$ws.Range("E31").Value = "  +3.30%  "
$ws.Range("D32").Value = "2.25"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("E33").Value = "  +1.06%  "
$ws.Range("D34").Value = "28.85"
$ws.Range("E34").Value = "  +0.08%  "
$ws.Range("D35").Value = "0.168"
$ws.Range("E35").Value = "  +18.16%  "
$ws.Range("D36").Value = "1.00"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "3.735.59"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").Value = "8.99"
$ws.Range("E38").Value = "  +0.68%  "
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").Value = "3.43"
$ws.Range("E40").Value = "  +8.22%  "
$ws.Range("E41").Value = "  +0.86%  "
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").Value = "0.999"
$ws.Range("E43").Value = "  -0.05%  "
$ws.Range("E45").Value = "  +0.40%  "
$ws.Range("E46").Value = "  +1.15%  "
$ws.Range("E47").Value = "  +4.25%  "
$ws.Range("D48").Value = "46.73"
$ws.Range("E48").Value = "  -0.84%  "
$ws.Range("D49").Value = "151.85"
$ws.Range("E49").Value = "  -0.19%  "
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("E51").Value = "  +0.40%  "
